$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 51, shifting existing rows 51-55 down to 52-56
$ws.Rows.Item(51).Insert()

$ws.Cells.Item(51, 1).Value = 9
$ws.Cells.Item(51, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(51, 3).Value = "Metropolitana"
$ws.Cells.Item(51, 4).Value = 44783
$ws.Cells.Item(51, 4).NumberFormat = $ws.Cells.Item(52, 4).NumberFormat
$ws.Cells.Item(51, 5).Value = 13
$ws.Cells.Item(51, 6).Value = 100112035
$ws.Cells.Item(51, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 35
$ws.Cells.Item(51, 11).Value = 20000
$ws.Cells.Item(51, 12).Value = 21000
$ws.Cells.Item(51, 13).Value = 20429
$ws.Cells.Item(51, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(51, 15).Value = "Hijuelas"
$ws.Cells.Item(51, 16).Value = 1362
$ws.Cells.Item(51, 17).Value = 15
$ws.Cells.Item(51, 18).Value = "Hortaliza"
